$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("tested")

# Fill in the Acc No. for Sheet1 row 2 (was blank/space placeholder)
$ws1.Range("A2").Value = "P0519961"

# Append the same record as a new row (row 4) on the "tested" sheet
$ws2.Cells.Item(4, 1).Value = "P0519961"
$ws2.Cells.Item(4, 2).Value = 13
$ws2.Cells.Item(4, 3).Value = "SGD"
$ws2.Cells.Item(4, 4).Value = "IFast"
$ws2.Cells.Item(4, 5).Value = "Wrap 0%"
$ws2.Cells.Item(4, 6).Value = "Top Up"
$ws2.Cells.Item(4, 7).Value = "Brown, Jarrad"
$ws2.Cells.Item(4, 8).Value = 0
$ws2.Cells.Item(4, 9).Value = 1
$ws2.Cells.Item(4, 10).Value = $ws1.Cells.Item(2, 10).Value()

# Match formatting/styles with the source row on Sheet1
$ws1.Range("A2:J2").Copy()
$ws2.Range("A4:J4").PasteSpecial(-4122)

# Update selection on the "tested" sheet to the newly added row, then
# restore Sheet1 as the active/visible tab (selecting on ws2 otherwise
# switches the active sheet)
$ws2.Range("A4").Select() | Out-Null
$ws1.Activate() | Out-Null
